$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the random data values in B2:D5 ---
$ws.Range("B2").Value = 0.24620924007686884
$ws.Range("C2").Value = 0.7744770598157588
$ws.Range("D2").Value = 0.30300638625243004

$ws.Range("B3").Value = 0.3630417644531141
$ws.Range("C3").Value = 0.2341433777172771
$ws.Range("D3").Value = 0.18752857485013474

$ws.Range("B4").Value = 0.6505125182365086
$ws.Range("C4").Value = 0.18576120723450062
$ws.Range("D4").Value = 0.22991321372747697

$ws.Range("B5").Value = 0.30227479120973777
$ws.Range("C5").Value = 0.2820220508202078
$ws.Range("D5").Value = 0.3588360230309111

# --- Widen the chart series ranges from column D to column E ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$chart.SeriesCollection().Item(1).Formula = '=SERIES("series1",Sheet1!$A$1:$E$1,Sheet1!$A$2:$E$2,1)'
$chart.SeriesCollection().Item(2).Formula = '=SERIES("series4",Sheet1!$A$1:$E$1,Sheet1!$A$3:$E$3,2)'
$chart.SeriesCollection().Item(3).Formula = '=SERIES(,Sheet1!$A$1:$E$1,Sheet1!$A$4:$E$4,3)'
$chart.SeriesCollection().Item(4).Formula = '=SERIES(,Sheet1!$A$1:$E$1,Sheet1!$A$5:$E$5,4)'
